# Apply the "Segunda actualización fichas tecnicas" edit.
#
# 1) Remove the plain paragraph that starts with
#    "El restaurante La Pescadería ha identificado..."
# 2) Add an underline to the (previously empty) paragraph that follows it,
#    and insert a new run - styled with character style "s1ppyq" and black
#    font color - containing the updated wording of that sentence.
# 3) Delete the now-unused "Normal (Web)" (NormalWeb) paragraph style.

$d = $word.ActiveDocument

# --- 1. Locate & delete the old paragraph -------------------------------
$oldText = "El restaurante La Pescadería ha identificado que el proceso de toma de pedidos y entrega de alimentos a los clientes puede resultar lento y poco eficiente durante horas pico debido al alto volumen de comensales y al proceso manual de toma de pedidos. Esto puede conducir a errores y retrasos, lo que afecta la satisfacción del cliente y puede tener un impacto negativo en la rentabilidad del restaurante."

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq $oldText) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}

# --- 2. Find the paragraph that used to be empty (now holds the new text) ---
$newParaText = "El restaurante La pescadería ha identificado que el proceso de toma de pedidos y entrega de alimentos a los clientes puede resultar lento y poco eficiente durante horas pico. Esto se debe a la gran cantidad de comensales que acuden al establecimiento y al proceso manual de toma de pedidos, que a menudo conduce a errores y retrasos. " + [string][char]0x200B

$destPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd("`r", "`a")
    if ($t -eq "") {
        $destPara = $p
        break
    }
}

if ($destPara -ne $null) {
    # Add the single underline to the paragraph mark / paragraph formatting.
    $destPara.Range.Font.Underline = 1

    # Insert the new run just before the paragraph mark.
    $insertRange = $destPara.Range
    $insertRange.Collapse(0)
    $insertRange.MoveEnd(1, -1)
    $insertRange.InsertBefore($newParaText)

    # Re-acquire the range for just the inserted text so formatting only
    # applies to the new run, not to the paragraph mark.
    $start = $destPara.Range.Start
    $end = $start + $newParaText.Length
    $runRange = $d.Range($start, $end)
    $runRange.Style = "s1ppyq"
    $runRange.Font.Color = 0
}

# --- 3. Remove the now unused "Normal (Web)" style -----------------------
foreach ($s in $d.Styles) {
    if ($s.NameLocal -eq "Normal (Web)") {
        $s.Delete()
        break
    }
}
